$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells G1, H1 with the same style as the existing header row (e.g. F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats

# Update the metric values in row 2
$ws.Range("B2").Value = 0.5086347285969581
$ws.Range("C2").Value = 0.9898720024061114
$ws.Range("D2").Value = 0.5812978092838879

# Update the model description text (now fits on a single line)
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(learning_rate=0.5))])"

# New Elapsed Time / CPU values for row 2
$ws.Range("G2").Value = 0.1228586025167412
$ws.Range("H2").Value = 0.991
